# "up truoc khi di hoc" - add the "bqt" (ban quan tri) table description
# as a new row at the bottom of Sheet1, and move the active selection
# down to the newly added description cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label in column A
$ws.Range("A10").Value = "bqt"

# New multi-colour description in column B, built the same way the other
# rows in the sheet describe their table's columns: the primary-key name
# in red, the plain attribute list in black, and the foreign-key list in
# blue.
$idPart = "banquantriId"
$attrPart = ",  firstName, lastName, gioiTinh, ngaySinh, soDienThoai, gmail, diaChi, "
$fkPart = "idCoSo , giangVienId,"

$ws.Range("B10").Value = $idPart + $attrPart + $fkPart

$idLen = $idPart.Length
$attrLen = $attrPart.Length
$fkLen = $fkPart.Length

# Red (FF0000) for the primary key
$ws.Range("B10").Characters(1, $idLen).Font.Color = 255
# Black (theme text colour) for the plain attributes
$ws.Range("B10").Characters($idLen + 1, $attrLen).Font.Color = 0
# Blue (accent colour) for the foreign keys
$ws.Range("B10").Characters($idLen + $attrLen + 1, $fkLen).Font.Color = 13998939

# Move the selection to B8, matching the author's final cursor position.
[void]$ws.Range("B8").Select()
